# Auto-generated from the cryptos.xlsx OOXML diff: refresh the
# Price (D) / Volume(1h) (E) columns for rows 2-51.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '53.950.34'
$ws.Range("E2").Value = '  +1.85%  '
$ws.Range("D3").Value = '2.236.93'
$ws.Range("E3").Value = '  +2.39%  '
$ws.Range("D4").Value = "'1.01"
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").Value = "'494.61"
$ws.Range("E5").Value = '  +4.06%  '
$ws.Range("D6").Value = "'127.37"
$ws.Range("E6").Value = '  +4.09%  '
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = '  -0.56%  '
$ws.Range("D8").Value = "'0.527"
$ws.Range("E8").Value = '  +3.08%  '
$ws.Range("D9").Value = '2.273.72'
$ws.Range("E9").Value = '  +3.92%  '
$ws.Range("D10").Value = "'0.0950"
$ws.Range("E10").Value = '  +5.50%  '
$ws.Range("D11").Value = "'0.152"
$ws.Range("E11").Value = '  +2.89%  '
$ws.Range("D12").Value = "'0.325"
$ws.Range("E12").Value = '  +5.63%  '
$ws.Range("D13").Value = "'4.64"
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("D14").Value = '2.645.26'
$ws.Range("E14").Value = '  +2.58%  '
$ws.Range("D15").Value = "'21.70"
$ws.Range("E15").Value = '  +4.76%  '
$ws.Range("D16").Value = '54.001.81'
$ws.Range("E16").Value = '  +1.99%  '
$ws.Range("E17").Value = '  +2.73%  '
$ws.Range("D18").Value = '2.255.75'
$ws.Range("E18").Value = '  +2.29%  '
$ws.Range("D19").Value = "'10.00"
$ws.Range("E19").Value = '  +6.47%  '
$ws.Range("E20").Value = '  +5.33%  '
$ws.Range("D21").Value = "'300.12"
$ws.Range("E21").Value = '  +2.51%  '
$ws.Range("E22").Value = '  +7.43%  '
$ws.Range("D23").Value = "'0.996"
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("E24").Value = '  -2.20%  '
$ws.Range("D25").Value = "'62.09"
$ws.Range("E25").Value = '  -0.61%  '
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = "'0.371"
$ws.Range("E27").Value = '  +2.95%  '
$ws.Range("D28").Value = '2.371.31'
$ws.Range("E28").Value = '  +3.17%  '
$ws.Range("E29").Value = '  +5.27%  '
$ws.Range("D30").Value = "'7.05"
$ws.Range("E30").Value = '  +2.97%  '
$ws.Range("D31").Value = "'167.87"
$ws.Range("E31").Value = '  +0.78%  '
$ws.Range("D32").Value = "'1.61"
$ws.Range("E32").Value = '  +3.90%  '
$ws.Range("D33").Value = '0.0₃0685'
$ws.Range("E33").Value = '  +5.12%  '
$ws.Range("D34").Value = "'5.86"
$ws.Range("E34").Value = '  +5.43%  '
$ws.Range("E35").Value = '  -0.14%  '
$ws.Range("D36").Value = "'0.991"
$ws.Range("E36").Value = '  -1.02%  '
$ws.Range("E37").Value = '  +4.85%  '
$ws.Range("D38").Value = "'17.70"
$ws.Range("E38").Value = '  +4.18%  '
$ws.Range("D39").Value = "'0.904"
$ws.Range("E39").Value = '  +12.60%  '
$ws.Range("E40").Value = '  +5.74%  '
$ws.Range("D41").Value = "'3.68"
$ws.Range("E41").Value = '  +5.93%  '
$ws.Range("D42").Value = "'35.58"
$ws.Range("E42").Value = '  +0.32%  '
$ws.Range("E43").Value = '  +5.38%  '
$ws.Range("E44").Value = '  +3.03%  '
$ws.Range("E45").Value = '  +4.88%  '
$ws.Range("D46").Value = "'125.91"
$ws.Range("E46").Value = '  +4.59%  '
$ws.Range("D47").Value = "'4.75"
$ws.Range("E47").Value = '  +4.17%  '
$ws.Range("E48").Value = '  +2.40%  '
$ws.Range("D49").Value = "'0.543"
$ws.Range("E49").Value = '  +3.82%  '
$ws.Range("D50").Value = "'236.75"
$ws.Range("E50").Value = '  +5.47%  '
$ws.Range("D51").Value = "'0.0483"
$ws.Range("E51").Value = '  +4.50%  '
